$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) and Volume(1h) (E) columns with the latest crypto data.
# Values in column D are stored as plain text in the workbook (General format).
# For values that look like plain numbers (e.g. "10.40", "2.00"), a leading
# apostrophe is used so Excel keeps them as text (preserving trailing zeros)
# instead of auto-converting them to numeric values.

$ws.Range("D2").Value = "30.723.52"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").Value = "1.686.61"
$ws.Range("E3").Value = "  +2.85%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'220.42"
$ws.Range("E5").Value = "  +2.19%  "

$ws.Range("D6").Value = "'0.525"
$ws.Range("E6").Value = "  +1.05%  "

$ws.Range("E7").Value = "  -0.30%  "

$ws.Range("D8").Value = "'30.44"
$ws.Range("E8").Value = "  +4.98%  "

$ws.Range("D9").Value = "'0.264"
$ws.Range("E9").Value = "  +1.85%  "

$ws.Range("D10").Value = "'0.0626"
$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").Value = "'0.0907"
$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("D13").Value = "'10.40"
$ws.Range("E13").Value = "  +11.28%  "

$ws.Range("D14").Value = "'0.620"
$ws.Range("E14").Value = "  +8.70%  "

$ws.Range("D15").Value = "1.684.96"
$ws.Range("E15").Value = "  +2.60%  "

$ws.Range("D16").Value = "'3.99"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").Value = "30.713.82"
$ws.Range("E17").Value = "  +2.34%  "

$ws.Range("D18").Value = "'66.35"
$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("D19").Value = "'246.34"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "0.0₃0715"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").Value = "  +1.77%  "

$ws.Range("D23").Value = "'4.29"
$ws.Range("E23").Value = "  +2.96%  "

$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'158.04"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").Value = "'15.88"
$ws.Range("E26").Value = "  +1.30%  "

$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("D28").Value = "'6.69"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("E31").Value = "  +0.86%  "

$ws.Range("D32").Value = "'3.48"
$ws.Range("E32").Value = "  +2.01%  "

$ws.Range("D33").Value = "1.511.87"
$ws.Range("E33").Value = "  +5.48%  "

$ws.Range("E34").Value = "  +3.10%  "

$ws.Range("E35").Value = "  +4.49%  "

$ws.Range("D36").Value = "'84.69"
$ws.Range("E36").Value = "  +10.12%  "

$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("E38").Value = "  +3.57%  "

$ws.Range("E39").Value = "  -4.92%  "

$ws.Range("D40").Value = "'0.586"
$ws.Range("E40").Value = "  +5.00%  "

$ws.Range("D41").Value = "'2.32"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").Value = "'0.838"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "'2.00"
$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D47").Value = "'52.13"
$ws.Range("E47").Value = "  -3.99%  "

$ws.Range("E48").Value = "  +1.92%  "

$ws.Range("D49").Value = "'5.43"
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("E50").Value = "  +5.45%  "

$ws.Range("E51").Value = "  +0.42%  "
